$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the mislabeled "number" data type to "integer" for the "A Number" row
$ws.Range("D5").Value = "integer"

# Add an explicit test row for a bad/unrecognized data type
$ws.Range("B6").Value = "Bad Type"
$ws.Range("C6").Value = "bad_type"
$ws.Range("D6").Value = "bad_type"

# Match the author's final selection/active cell state
$ws.Range("B7").Select()
